$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 113 - this shifts the existing rows 113-120
# down to 114-121, matching the rest of the diff automatically.
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new data point.
$ws.Range("A113").Value = 3
$ws.Range("B113").Value = "Femacal de La Calera"
$ws.Range("C113").Value = "Coquimbo"
$ws.Range("D113").Value = 44578
$ws.Range("E113").Value = 5
$ws.Range("F113").Value = 100112052
$ws.Range("G113").Value = "Albahaca"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 130
$ws.Range("K113").Value = 4500
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 4769
$ws.Range("N113").Value = "$/docena de matas"
$ws.Range("O113").Value = "Provincia de Quillota"
$ws.Range("P113").Value = 795
$ws.Range("Q113").Value = 6
$ws.Range("R113").Value = "Hortaliza"
